$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.926.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.546.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.97%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.97%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.578"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.85%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.548"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0832"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.91%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.77"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.78%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.115"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.72%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.942.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.509.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.875"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.74%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.959.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.69%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0998"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.55%  "

$ws.Range("E21").Value = "  +1.05%  "

$ws.Range("E22").Value = "  +0.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "255.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.33%  "

$ws.Range("E24").Value = "  +1.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.53%  "

$ws.Range("E30").Value = "  -1.13%  "

$ws.Range("E31").Value = "  +4.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.24%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.69"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +17.03%  "

$ws.Range("E34").Value = "  -1.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0802"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.23%  "

$ws.Range("E36").Value = "  -2.17%  "

$ws.Range("E37").Value = "  -4.51%  "

$ws.Range("E38").Value = "  +2.57%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.21%  "

$ws.Range("E40").Value = "  +0.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.68%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.53%  "

$ws.Range("E43").Value = "  +28.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.101.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.48%  "

$ws.Range("E45").Value = "  -1.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.47%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.798.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.192"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.54%  "
